$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.738.80"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "1.751.05"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "'235.94"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").Value = "'0.5046"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").Value = "'41.05"
$ws.Range("E8").Value = "  -4.67%  "
$ws.Range("E9").Value = "  +6.93%  "
$ws.Range("D10").Value = "'0.06242"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").Value = "1.748.67"
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("D12").Value = "'0.06921"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "'15.44"
$ws.Range("E13").Value = "  +4.49%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.5980"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'4.471"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "'76.72"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "25.752.18"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000006800"
$ws.Range("E20").Value = "  +16.18%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'11.59"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").Value = "1.971.84"
$ws.Range("E22").Value = "  -4.57%  "
$ws.Range("D23").Value = "'4.067"
$ws.Range("E23").Value = "  +3.27%  "
$ws.Range("D24").Value = "'8.230"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").Value = "'5.172"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "'136.71"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("D27").Value = "'1.445"
$ws.Range("E27").Value = "  +16.08%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.806"
$ws.Range("E28").Value = "  -3.28%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'14.94"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("D30").Value = "'102.12"
$ws.Range("E30").Value = "  +3.65%  "
$ws.Range("D31").Value = "'0.08161"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").Value = "'3.662"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "'3.415"
$ws.Range("E33").Value = "  +8.66%  "
$ws.Range("D34").Value = "'0.04473"
$ws.Range("E34").Value = "  +4.80%  "
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "'2.658"
$ws.Range("E36").Value = "  -4.02%  "
$ws.Range("D37").Value = "'0.9905"
$ws.Range("E37").Value = "  -5.94%  "
$ws.Range("D38").Value = "'0.6041"
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("D39").Value = "'2.686"
$ws.Range("E39").Value = "  -8.09%  "
$ws.Range("D40").Value = "'0.01548"
$ws.Range("E40").Value = "  +6.00%  "
$ws.Range("D41").Value = "'1.934"
$ws.Range("E41").Value = "  -7.28%  "
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'102.66"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("D44").Value = "'0.3785"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").Value = "'0.7361"
$ws.Range("E45").Value = "  -11.22%  "
$ws.Range("E46").Value = "  -6.05%  "
$ws.Range("D47").Value = "'0.05470"
$ws.Range("E47").Value = "  +4.07%  "
$ws.Range("D48").Value = "'0.1094"
$ws.Range("E48").Value = "  +6.32%  "
$ws.Range("D49").Value = "'5.908"
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("D50").Value = "'7.621"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").Value = "'29.71"
$ws.Range("E51").Value = "  +0.01%  "
